# Auto-generated Excel COM-interop script
# Applies numeric value updates to H:N columns across multiple sheets
# as described by the commit diff (scheduled runner profit recompute).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 73642.734
$ws.Range("I9").Value = 110171.6
$ws.Range("J9").Value = 585
$ws.Range("K9").Value = 110171.6
$ws.Range("L9").Value = 585
$ws.Range("M9").Value = -110002.6
$ws.Range("N9").Value = -923
$ws.Range("H12").Value = 530.1539
$ws.Range("I12").Value = 124
$ws.Range("J12").Value = 1180
$ws.Range("K12").Value = 124
$ws.Range("L12").Value = 1180
$ws.Range("M12").Value = 46
$ws.Range("N12").Value = -1520
$ws.Range("H33").Value = 283.4
$ws.Range("I33").Value = 215.5
$ws.Range("K33").Value = 215.5
$ws.Range("M33").Value = 13.5
$ws.Range("H40").Value = 3300.2
$ws.Range("I40").Value = 2875.25
$ws.Range("K40").Value = 2875.25
$ws.Range("M40").Value = -2700.25
$ws.Range("H51").Value = 24999.5
$ws.Range("J51").Value = 9999
$ws.Range("L51").Value = 9999
$ws.Range("N51").Value = -10967
$ws.Range("H58").Value = 2036
$ws.Range("I58").Value = 381.33334
$ws.Range("K58").Value = 1144.00002
$ws.Range("M58").Value = -994.0000199999999
$ws.Range("H64").Value = 10818
$ws.Range("J64").Value = 12002.5
$ws.Range("L64").Value = 12002.5
$ws.Range("N64").Value = -12498.5
$ws.Range("H67").Value = 10818
$ws.Range("J67").Value = 12002.5
$ws.Range("L67").Value = 12002.5
$ws.Range("N67").Value = -13718.5
$ws.Range("H80").Value = 618
$ws.Range("I80").Value = 398
$ws.Range("K80").Value = 1194
$ws.Range("M80").Value = -196
$ws.Range("H83").Value = 618
$ws.Range("I83").Value = 398
$ws.Range("K83").Value = 3582
$ws.Range("M83").Value = 1410
$ws.Range("H96").Value = 1211.9231
$ws.Range("I96").Value = 1261.4286
$ws.Range("J96").Value = 1154.1666
$ws.Range("K96").Value = 3784.2858
$ws.Range("L96").Value = 3462.4998
$ws.Range("M96").Value = -2411.2858
$ws.Range("N96").Value = -6208.4998
$ws.Range("H125").Value = 2525.5386
$ws.Range("I125").Value = 1074.8889
$ws.Range("K125").Value = 9674.000099999999
$ws.Range("M125").Value = -7214.000099999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4598.827
$ws.Range("I32").Value = 3645.9387
$ws.Range("K32").Value = 3645.9387
$ws.Range("M32").Value = -3358.9387
$ws.Range("H45").Value = 3569.318
$ws.Range("I45").Value = 2599.4285
$ws.Range("K45").Value = 2599.4285
$ws.Range("M45").Value = -2222.4285
$ws.Range("H61").Value = 6757
$ws.Range("I61").Value = 3270.611
$ws.Range("K61").Value = 3270.611
$ws.Range("M61").Value = -3058.611
$ws.Range("H74").Value = 2040.3077
$ws.Range("I74").Value = 1391.1111
$ws.Range("K74").Value = 1391.1111
$ws.Range("M74").Value = -517.1111000000001
$ws.Range("H77").Value = 2040.3077
$ws.Range("I77").Value = 1391.1111
$ws.Range("K77").Value = 6955.5555
$ws.Range("M77").Value = -2587.5555
$ws.Range("H97").Value = 1955.8572
$ws.Range("I97").Value = 622.9259
$ws.Range("K97").Value = 622.9259
$ws.Range("M97").Value = -126.9259
$ws.Range("H136").Value = 6757
$ws.Range("I136").Value = 3270.611
$ws.Range("K136").Value = 9811.832999999999
$ws.Range("M136").Value = -7261.832999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 274722
$ws.Range("J70").Value = 274722
$ws.Range("L70").Value = 274722
$ws.Range("N70").Value = -275308
$ws.Range("H73").Value = 274722
$ws.Range("J73").Value = 274722
$ws.Range("L73").Value = 274722
$ws.Range("N73").Value = -276750
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1922.4584
$ws.Range("I58").Value = 1457.7222
$ws.Range("K58").Value = 1457.7222
$ws.Range("M58").Value = -1254.7222
$ws.Range("H105").Value = 1441.3334
$ws.Range("I105").Value = 1183
$ws.Range("J105").Value = 1699.6666
$ws.Range("K105").Value = 1183
$ws.Range("L105").Value = 1699.6666
$ws.Range("M105").Value = 564
$ws.Range("N105").Value = -5193.6666
$ws.Range("H132").Value = 3053.75
$ws.Range("I132").Value = 2723.5
$ws.Range("K132").Value = 8170.5
$ws.Range("M132").Value = -5640.5
$ws.Range("H134").Value = 2240
$ws.Range("I134").Value = 1481.8334
$ws.Range("J134").Value = 4189.5713
$ws.Range("K134").Value = 4445.5002
$ws.Range("L134").Value = 12568.7139
$ws.Range("M134").Value = -1910.5002
$ws.Range("N134").Value = -17638.7139
$ws.Range("H136").Value = 1922.4584
$ws.Range("I136").Value = 1457.7222
$ws.Range("K136").Value = 4373.1666
$ws.Range("M136").Value = -1823.1666
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 730.2222
$ws.Range("I26").Value = 259.2
$ws.Range("J26").Value = 911.38464
$ws.Range("K26").Value = 777.5999999999999
$ws.Range("L26").Value = 2734.15392
$ws.Range("M26").Value = -489.5999999999999
$ws.Range("N26").Value = -3310.15392
$ws.Range("H131").Value = 24692828
$ws.Range("I131").Value = 15873988
$ws.Range("J131").Value = 30304818
$ws.Range("K131").Value = 47621964
$ws.Range("L131").Value = 90914454
$ws.Range("M131").Value = -47616924
$ws.Range("N131").Value = -90924534
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 13444.5
$ws.Range("I27").Value = 13444
$ws.Range("K27").Value = 13444
$ws.Range("M27").Value = -13278
$ws.Range("H107").Value = 1087.5454
$ws.Range("J107").Value = 2500
$ws.Range("L107").Value = 2500
$ws.Range("N107").Value = -6340
$ws.Range("H132").Value = 2534.4146
$ws.Range("J132").Value = 4290.4546
$ws.Range("L132").Value = 12871.3638
$ws.Range("N132").Value = -17931.3638
$ws.Range("H136").Value = 27068.875
$ws.Range("J136").Value = 27068.875
$ws.Range("L136").Value = 81206.625
$ws.Range("N136").Value = -86306.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1452.3334
$ws.Range("I22").Value = 1310.5555
$ws.Range("K22").Value = 1310.5555
$ws.Range("M22").Value = -1015.5555
$ws.Range("H27").Value = 1452.3334
$ws.Range("I27").Value = 1310.5555
$ws.Range("K27").Value = 1310.5555
$ws.Range("M27").Value = -1203.5555
$ws.Range("H46").Value = 2030.4375
$ws.Range("J46").Value = 2509.6667
$ws.Range("L46").Value = 2509.6667
$ws.Range("N46").Value = -2885.6667
$ws.Range("H55").Value = 925.52
$ws.Range("I55").Value = 303.54544
$ws.Range("J55").Value = 1414.2142
$ws.Range("K55").Value = 303.54544
$ws.Range("L55").Value = 1414.2142
$ws.Range("M55").Value = -130.54544
$ws.Range("N55").Value = -1760.2142
$ws.Range("H68").Value = 2440.5
$ws.Range("I68").Value = 2433.889
$ws.Range("K68").Value = 2433.889
$ws.Range("M68").Value = -1684.889
$ws.Range("H71").Value = 2440.5
$ws.Range("I71").Value = 2433.889
$ws.Range("K71").Value = 12169.445
$ws.Range("M71").Value = -8425.445
$ws.Range("H82").Value = 7858.1665
$ws.Range("I82").Value = 1774
$ws.Range("K82").Value = 1774
$ws.Range("M82").Value = -1413
$ws.Range("H85").Value = 7858.1665
$ws.Range("I85").Value = 1774
$ws.Range("K85").Value = 1774
$ws.Range("M85").Value = -526
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 15000
$ws.Range("M113").Value = -12830
$ws.Range("H132").Value = 3211.7693
$ws.Range("I132").Value = 3417.111
$ws.Range("K132").Value = 10251.333
$ws.Range("M132").Value = -7721.332999999999
$ws.Range("H136").Value = 9579.739
$ws.Range("I136").Value = 7512.684
$ws.Range("K136").Value = 22538.052
$ws.Range("M136").Value = -19988.052
